$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ydh_alert_list")

$ws.Rows.Item(83).Insert()
$rng = $ws.Range("A83:E83")
$rng.Font.Name = "Calibri"

$ws.Range("A83").Value = "WN-CoV"
$ws.Range("B83").Value = "WN-Cov swab taken"
$ws.Range("E83").Value = "Patient swabbed for coronavirus / Covid-19"

$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:E84"))
